# Update "paises" workbook: refresh COVID country stats and reorder Japon
# in the country list (shared-string order swap), plus bump the "datos
# actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Agosto de 2020 a las 01:48"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 4762954
$ws.Range("C4").Value = 57065
$ws.Range("D4").Value = 2355848
$ws.Range("E4").Value = 2249251
$ws.Range("G4").Value = 1080
$ws.Range("H4").Value = 157855

# --- Row 5: Brasil ---
$ws.Range("B5").Value = 2708876
$ws.Range("C5").Value = 42578
$ws.Range("E5").Value = 731209
$ws.Range("G5").Value = 1048
$ws.Range("H5").Value = 93616

# --- Row 50: Nigeria ---
$ws.Range("B50").Value = 43537
$ws.Range("C50").Value = 386
$ws.Range("D50").Value = 20087
$ws.Range("E50").Value = 22567
$ws.Range("G50").Value = 4
$ws.Range("H50").Value = 883

# --- Rows 56-58: Japon moves above Ghana/Suiza in the country ordering.
# Row 56 now shows Japon (with fresh data), row 57 now shows Ghana (its
# old, unchanged data shifted down one row), row 58 now shows Suiza
# (its old, unchanged data shifted down one row). Azerbaiyan (row 59)
# is untouched.
$ws.Range("A56").Value = "Japon"
$ws.Range("B56").Value = 35836
$ws.Range("C56").Value = 1464
$ws.Range("D56").Value = 25506
$ws.Range("E56").Value = 9319
$ws.Range("G56").Value = 5
$ws.Range("H56").Value = 1011

$ws.Range("A57").Value = "Ghana"
$ws.Range("B57").Value = 35501
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 32096
$ws.Range("E57").Value = 3223
$ws.Range("H57").Value = 182

$ws.Range("A58").Value = "Suiza"
$ws.Range("B58").Value = 35412
$ws.Range("C58").Value = 180
$ws.Range("D58").Value = 31300
$ws.Range("E58").Value = 2131
$ws.Range("H58").Value = 1981

# --- Row 75: Chequia ---
$ws.Range("B75").Value = 16699
$ws.Range("C75").Value = 291
$ws.Range("D75").Value = 11587
$ws.Range("E75").Value = 4729
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 383

# --- Row 86: Noruega ---
$ws.Range("B86").Value = 9253
$ws.Range("C86").Value = 13
$ws.Range("E86").Value = 246

# --- Row 95: Luxemburgo ---
$ws.Range("B95").Value = 6793
$ws.Range("C95").Value = 98
$ws.Range("E95").Value = 1485
$ws.Range("G95").Value = 2
$ws.Range("H95").Value = 116

# --- Row 130: Guinea-Bisau ---
$ws.Range("E130").Value = 1151
$ws.Range("G130").Value = 1
$ws.Range("H130").Value = 27

# --- Row 139: Uruguay ---
$ws.Range("B139").Value = 1278
$ws.Range("C139").Value = 14
$ws.Range("D139").Value = 1004
$ws.Range("E139").Value = 239

# --- Row 147: Niger ---
$ws.Range("B147").Value = 1136
$ws.Range("C147").Value = 2
$ws.Range("E147").Value = 39

# --- Row 161: Vietnam ---
$ws.Range("B161").Value = 590
$ws.Range("C161").Value = 44
$ws.Range("E161").Value = 214
